# Adani-BUWise-Websites.xlsx edit script
# - Update RMRW row (row 2) website list: drop "uat-s.crfindia.org" line
# - Update Airports row (row 5) website list: swap mangaluru URL format
# - Add a new row 19 for "DRP" business unit with its website + comment
# - Expand Table1 to cover the new row
# - Add a hyperlink for the new DRP website cell
# - Update the active selection / top-left cell to match final view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (RMRW): remove the "uat-s.crfindia.org" line from the Websites cell
$ws.Range("C2").Value = "www.bprpl.in`nwww.avrpl.com`nwww.bkrpl.com`nwww.kkrpl.com`nwww.mrrpl.com`nwww.nprpl.com`nwww.skrpl.com`nwww.vbppl.in`nwww.ksrpl.co.in"
$ws.Rows.Item(2).RowHeight = 145.8

# --- Row 5 (Airports): update the Mangaluru airport URL
$ws.Range("C5").Value = "www.adaniairports.com`nsvpia-ahmedabad.adaniairports.com`nwww.adani.com/ccsia-lucknow-airport`nmangaluru.adaniairports.com/`ncsmia-mumbai.adaniairports.com`nwww.adani.com/jaipur-airport`nwww.adani.com/lgbia-guwahati-airport`nthiruvananthapuram.adaniairports.com`nnavimumbai.adaniairports.com"

# --- New row 19: DRP business unit
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "DRP"
$ws.Range("C19").Value = "https://drpmumbai.maharashtra.gov.in/"
$ws.Range("D19").Value = "Need MVC webapp only"
$ws.Rows.Item(19).RowHeight = 16.2

# Hyperlink for the new DRP website cell
$ws.Hyperlinks.Add($ws.Range("C19"), "https://drpmumbai.maharashtra.gov.in/")

# Copy formatting from row 18 last, so the copied cell styles (incl. the
# hyperlink look for C19) win over whatever Hyperlinks.Add touched
$ws.Range("A18:D18").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)

# --- Expand the table to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D19"))

# --- Update the view: selection on C5, scrolled so row 4 is at top
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
